# Word COM-interop script implementing:
#   1. Append " (For FSET version 1.1 and above)" to the "Instructions" heading.
#   2. Add a new paragraph after the "...please find the SDK online and
#      download it." paragraph, announcing that the AreaKML.kml method has
#      been disabled.
#
# (The source diff also shows many w:proofErr "gramStart"/"gramEnd" wrapper
#  splits and one w:lastRenderedPageBreak marker; those are cosmetic
#  grammar-checker / page-layout rendering artifacts that carry no textual
#  content and cannot be produced through the Word object model, so they are
#  intentionally not reproduced here.)

$d = $word.ActiveDocument

# --- 1. Heading: "Instructions" -> "Instructions (For FSET version 1.1 and above)"
$d.Content.Find.Execute(
    "Instructions",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Instructions (For FSET version 1.1 and above)", 2) | Out-Null

# --- 2. Insert the new "Warning" paragraph after the addendum paragraph.
$anchor = $d.Content
$anchor.Find.Execute(
    "please find the SDK online and download it.",
    $true, $false, $false, $false, $false,
    $true, 0, $false, "", 0) | Out-Null

$endOfPara1 = $anchor.End

# Insert the first (blank spacer) paragraph right after the addendum paragraph.
$ip = $anchor.Duplicate
$ip.Collapse(0)
$ip.InsertParagraphAfter()

$newPara1Start = $endOfPara1 + 1

# Insert the second (content) paragraph right after the blank spacer.
$ip2 = $d.Range($newPara1Start, $newPara1Start)
$ip2.InsertParagraphAfter()

$newPara2Start = $newPara1Start + 1
$contentRange = $d.Range($newPara2Start, $newPara2Start)
$contentRange.InsertAfter("Warning " + [char]0x2013 + " the AreaKML.kml method of creating water masks inside FSET has been disabled for the time being since FSET automatically creates water masked sceneries. I might re-enable support for custom AreaKML files if many people request it or I see the need to do it.")
